# Insert a new data row at row 32 (pushing the existing rows 32-64 down to 33-65)
# and populate it with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(32).Insert()

$ws.Range("A32").Value = 10
$ws.Range("B32").Value = "Vega Modelo de Temuco"
$ws.Range("C32").Value = "La Araucanía"
$ws.Range("D32").Value = 44574
$ws.Range("E32").Value = 9
$ws.Range("F32").Value = 100112026
$ws.Range("G32").Value = "Haba"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 280
$ws.Range("K32").Value = 15000
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = 15000
$ws.Range("N32").Value = "`$/saco 25 kilos"
$ws.Range("O32").Value = "Región de La Araucanía"
$ws.Range("P32").Value = 600
$ws.Range("Q32").Value = 25
$ws.Range("R32").Value = "Hortaliza"
